$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 141, shifting existing rows 141-252 down to 142-253.
$ws.Rows.Item(141).Insert()

# Populate the newly inserted row 141 with the new record's data.
$ws.Range("A141").Value = 4
$ws.Range("B141").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C141").Value = "Los Lagos"
$ws.Range("D141").Value = 44673
$ws.Range("E141").Value = 10
$ws.Range("F141").Value = 100112003
$ws.Range("G141").Value = "Ajo"
$ws.Range("H141").Value = "Chino"
$ws.Range("I141").Value = "Primera"
$ws.Range("J141").Value = 60
$ws.Range("K141").Value = 19000
$ws.Range("L141").Value = 20000
$ws.Range("M141").Value = 19500
$ws.Range("N141").Value = "$/caja 10 kilos"
$ws.Range("O141").Value = "China"
$ws.Range("P141").Value = 1950
$ws.Range("Q141").Value = 10
$ws.Range("R141").Value = "Hortaliza"
